# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-10-30 08:19:08
#
# This session-analysis sheet is regenerated each sync: the "Recorded By"
# audit-trail lists get re-serialized (first contributor rotates to the
# end of the comma list), three B2D/B2E/B2F sessions that were still
# "Pending" have now elapsed and flip to "Not Recorded", and the Class /
# Group statistics panel is recomputed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Reorder "Recorded By" (G column) audit lists: rotate the first
#        contributor in the comma-separated list to the end. ---
$gRows = @(2,3,4,5,6,7,8,10,11,12,13,14,15,17,18,19,20,21,22,24,29,30,31,32,33,34,35,37,38,39,40,41,42,44,45,46,47,48,49,51,56,57,58,59,60,61,62,64,65,66,67,68,69,71,72,73,74,75,76,78,83,84,85,86,87,88,89,90,93,95,96,97,99,102,104,109,110,111,112,113,114,115,116,119,121,122,123,125,128,130,135,136,137,138,139,140,141,142,145,147,148,149,151,154,156)

foreach ($r in $gRows) {
  $cell = $ws.Cells.Item($r, 7)
  $v = $cell.Value2
  $parts = $v -split ", "
  $rot = $parts[1..($parts.Length - 1)] + $parts[0]
  $cell.Value = $rot -join ", "
}

# --- 2. The 30/10/2025 sessions for B2D (row 105), B2E (row 131) and
#        B2F (row 157) are now in the past but were never recorded:
#        flip their status from "Pending" (yellow) to "Not Recorded"
#        (red/pink) and update the cell color to match. ---
$statusRows = @(105, 131, 157)
foreach ($r in $statusRows) {
  $rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 9))
  $rowRange.Interior.Color = 12695295   # RGB(255,182,193) "FFB6C1" - Not Recorded / red
  $rowRange.Font.Color = 0              # black text
  $ws.Cells.Item($r, 9).Value = "Not Recorded"
}

# --- 3. Widen the Status column (I) so "Not Recorded" fits. ---
$ws.Columns.Item(9).ColumnWidth = 13.1667

# --- 4. Recompute the summary statistics affected by the 3 status flips:
#        +3 Missing Sessions / -3 Pending Sessions overall, and +1 Missing
#        / -1 Pending for each of the B2D/B2E/B2F group rows. ---
$ws.Range("L7").Value = 3    # Missing Sessions (overall class statistics)
$ws.Range("L8").Value = 18   # Pending Sessions (overall class statistics)

$ws.Range("P18").Value = 1   # B2D - Missing
$ws.Range("Q18").Value = 3   # B2D - Pending
$ws.Range("P19").Value = 1   # B2E - Missing
$ws.Range("Q19").Value = 3   # B2E - Pending
$ws.Range("P20").Value = 1   # B2F - Missing
$ws.Range("Q20").Value = 3   # B2F - Pending

Write-Output "Applied sync edit: reordered recorder lists, flipped 3 sessions to Not Recorded, refreshed stats."
